# Remove the temporary "Sheet1" scratch sheet and promote "Sheet2" to be
# the single, named "Product backlog" sheet. Also backfill the Sprint
# Number column (A) for a handful of rows that were missing it.

$wb = $excel.ActiveWorkbook

# Suppress the "this sheet may contain data" confirmation dialog so the
# delete below runs unattended, then restore the normal alert behavior.
$excel.DisplayAlerts = $false

# Delete the old temp sheet (Sheet1).
$wsOld = $wb.Worksheets.Item("Sheet1")
$wsOld.Delete()

$excel.DisplayAlerts = $true

# Rename the remaining sheet and make it active.
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "Product backlog"
$ws.Activate()

# Restore the 100% zoom level for the view.
$win = $wb.Windows.Item(1)
$win.Zoom = 100

# Fill in missing Sprint Number (column A) values for rows 2, 3, 5, 6, 8.
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 1
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 1
$ws.Range("A8").Value = 1

# Update the current selection to match the saved view (F14).
$ws.Range("F14").Select()
